$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.150.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.994.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.01%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.60%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.984.09"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("E10").Value = "  +7.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.17"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.44%  "
$ws.Range("E12").Value = "  +4.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.86%  "
$ws.Range("E15").Value = "  +2.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.489.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.994.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "59.135.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.97%  "
$ws.Range("E22").Value = "  +6.78%  "
$ws.Range("E23").Value = "  +3.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.55%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.43%  "
$ws.Range("E29").Value = "  +3.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0987"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0775"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +21.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.992"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.85%  "
$ws.Range("E37").Value = "  +2.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "401.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.778.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.53%  "
$ws.Range("E43").Value = "  +4.00%  "
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("E45").Value = "  +10.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("E48").Value = "  +2.33%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +22.70%  "
$ws.Range("B50").Value = "Fetch.AI"
$ws.Range("C50").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.92%  "
